{"js": "// Applies the cover-letter content update (AllSaints/Junior Front End Developer\n// -> MyEdSpace/Graduate Engineer) described by the diff.\n//\n// Order matters: the two long body paragraphs and the \"To the hiring team...\"\n// line contain \"AllSaints\" / \"Junior Front End Developer\" as substrings, so\n// those full-paragraph replacements run FIRST. Only afterwards do we replace\n// the short standalone phrases (heading subtitle, signature name), once they\n// are no longer ambiguous with leftover substring matches.\nconst replacements = [\n  [\n    \"To the hiring team at AllSaints\",\n    \"To the hiring team at MyEdSpace\",\n  ],\n  [\n    \"I am excited to apply for the Junior Front End Developer position at AllSaints. The role aligns perfectly with my skills and aspirations, especially in making customers feel cool and confident, a field that strongly interests me. AllSaints' focus on developing and optimizing website content pages resonates with my passion - having built a full-stack food ordering platform where I significantly improved customer experience through front-end optimizations, and I am eager to contribute while growing with your team.\",\n    \"I am excited to apply for the Graduate Engineer position at MyEdSpace. The role aligns perfectly with my skills and aspirations, especially in transforming education globally with technology, a field that strongly interests me. MyEdSpace's focus on building impactful features and collaborating on high-priority projects resonates with my passion - as a Full Stack Engineer, I have built a full-stack food ordering platform that improved customer experience and reduced communication overhead by 30%. I am eager to contribute while growing with your team.\",\n  ],\n  [\n    \"I am a Full Stack Engineer who recently built a full-stack food ordering platform with real-time order processing. This experience strengthened my proficiency in HTML, CSS, JavaScript, and responsive design and deepened my passion for solving practical challenges. A specific achievement from my previous experience that I believe can add value to the Junior Front End Developer position at AllSaints includes:\",\n    \"I am a Full Stack Engineer who recently built a full-stack food ordering platform with real-time processing. This experience strengthened my experience with React.js, Node.js, and AWS and deepened my passion for solving practical challenges. A specific achievement from my previous experience that I believe can add value to the Graduate Engineer position at MyEdSpace includes:\",\n  ],\n  [\n    \"Improved customer experience by 10% through internationalization.\",\n    \"Improved customer experience and reduced communication overhead by 30%.\",\n  ],\n  [\n    \"Experience in optimizing front-end features for enhanced performance.\",\n    \"Optimized SQL queries for backend performance.\",\n  ],\n  [\n    \"Real-time features significantly enhance user engagement.\",\n    \"Data-driven decisions can significantly enhance user engagement.\",\n  ],\n  [\n    \"My unique background as a Full Stack Engineer | Marketing Content Management Platform has provided me with experience in designing a multi-version content management system that enhances team collaboration, which I believe can also contribute to driving the company\\u2019s success in achieving the company's goal.\",\n    \"My unique background as a Full Stack Engineer on a Marketing Content Management Platform has provided me with experience in designing and building collaborative content management systems, integrated with AI for improved efficiency, which I believe can also contribute to driving the company\\u2019s success in achieving the company's goal.\",\n  ],\n  // Short standalone phrases last, once the longer paragraphs that used to\n  // contain them as substrings have already been rewritten.\n  [\"Junior Front End Developer\", \"Graduate Engineer\"],\n  [\"Amy Han Hsun Shi\", \"Amy Han Hsun Shiha\"],\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of replacements) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    throw new Error(\"No match found for: \" + oldText);\n  }\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Applies the cover-letter content update (AllSaints/Junior Front End\n# Developer -> MyEdSpace/Graduate Engineer) described by the diff, using\n# Word COM interop (Find/Replace across the whole document).\n\n$d = $word.ActiveDocument\n\n# wdReplace.wdReplaceAll = 2, wdFindWrap.wdFindContinue = 1\n$wdReplaceAll = 2\n$wdFindContinue = 1\n\nfunction Replace-AllText($document, [string]$findText, [string]$replaceText, [bool]$matchCase) {\n    $rng = $document.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Replacement.ClearFormatting()\n    $ok = $rng.Find.Execute(\n        $findText,       # FindText\n        $matchCase,      # MatchCase\n        $false,          # MatchWholeWord\n        $false,          # MatchWildcards\n        $false,          # MatchSoundsLike\n        $false,          # MatchAllWordForms\n        $true,           # Forward\n        $wdFindContinue, # Wrap\n        $false,          # Format\n        $replaceText,    # ReplaceWith\n        $wdReplaceAll    # Replace\n    )\n    if (-not $ok) {\n        throw \"No match found for: $findText\"\n    }\n}\n\n# Order matters: the two long body paragraphs and the \"To the hiring team...\"\n# line contain \"AllSaints\" / \"Junior Front End Developer\" as substrings, so\n# those full-paragraph replacements run FIRST. Only afterwards do we replace\n# the short standalone phrases (heading subtitle, signature name), once they\n# are no longer ambiguous with leftover substring matches.\n\nReplace-AllText $d \"To the hiring team at AllSaints\" \"To the hiring team at MyEdSpace\" $true\n\n$find1 = \"I am excited to apply for the Junior Front End Developer position at AllSaints. The role aligns perfectly with my skills and aspirations, especially in making customers feel cool and confident, a field that strongly interests me. AllSaints' focus on developing and optimizing website content pages resonates with my passion - having built a full-stack food ordering platform where I significantly improved customer experience through front-end optimizations, and I am eager to contribute while growing with your team.\"\n$repl1 = \"I am excited to apply for the Graduate Engineer position at MyEdSpace. The role aligns perfectly with my skills and aspirations, especially in transforming education globally with technology, a field that strongly interests me. MyEdSpace's focus on building impactful features and collaborating on high-priority projects resonates with my passion - as a Full Stack Engineer, I have built a full-stack food ordering platform that improved customer experience and reduced communication overhead by 30%. I am eager to contribute while growing with your team.\"\nReplace-AllText $d $find1 $repl1 $true\n\n$find2 = \"I am a Full Stack Engineer who recently built a full-stack food ordering platform with real-time order processing. This experience strengthened my proficiency in HTML, CSS, JavaScript, and responsive design and deepened my passion for solving practical challenges. A specific achievement from my previous experience that I believe can add value to the Junior Front End Developer position at AllSaints includes:\"\n$repl2 = \"I am a Full Stack Engineer who recently built a full-stack food ordering platform with real-time processing. This experience strengthened my experience with React.js, Node.js, and AWS and deepened my passion for solving practical challenges. A specific achievement from my previous experience that I believe can add value to the Graduate Engineer position at MyEdSpace includes:\"\nReplace-AllText $d $find2 $repl2 $true\n\nReplace-AllText $d \"Improved customer experience by 10% through internationalization.\" \"Improved customer experience and reduced communication overhead by 30%.\" $true\n\nReplace-AllText $d \"Experience in optimizing front-end features for enhanced performance.\" \"Optimized SQL queries for backend performance.\" $true\n\nReplace-AllText $d \"Real-time features significantly enhance user engagement.\" \"Data-driven decisions can significantly enhance user engagement.\" $true\n\n$find3 = \"My unique background as a Full Stack Engineer | Marketing Content Management Platform has provided me with experience in designing a multi-version content management system that enhances team collaboration, which I believe can also contribute to driving the company\" + [char]0x2019 + \"s success in achieving the company's goal.\"\n$repl3 = \"My unique background as a Full Stack Engineer on a Marketing Content Management Platform has provided me with experience in designing and building collaborative content management systems, integrated with AI for improved efficiency, which I believe can also contribute to driving the company\" + [char]0x2019 + \"s success in achieving the company's goal.\"\nReplace-AllText $d $find3 $repl3 $true\n\n# Short standalone phrases last, once the longer paragraphs that used to\n# contain them as substrings have already been rewritten.\n#\n# \"Junior Front End Developer\" lives in the subtitle heading under a\n# w:caps (all-caps display) run, so it must be matched case-insensitively\n# (MatchCase=$false) -- the displayed text is \"JUNIOR FRONT END DEVELOPER\"\n# even though the underlying run text is mixed case. The underlying stored\n# text, and hence the replacement text we supply, stays correctly cased.\nReplace-AllText $d \"Junior Front End Developer\" \"Graduate Engineer\" $false\nReplace-AllText $d \"Amy Han Hsun Shi\" \"Amy Han Hsun Shiha\" $true\n"}
